$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clean up header labels (drop trailing spaces) ---
$ws.Range("C1").Value = "Confirmados"
$ws.Range("D1").Value = "Negativos"
$ws.Range("E1").Value = "Sospechosos"
$ws.Range("F1").Value = "Defunciones"

# --- 2. Tidy floating point rounding noise in existing rows ---
$ws.Range("G3").Value = 34.59
$ws.Range("G4").Value = 34.45

# --- 3. Add the new row (June 5th) before re-styling B5 so we can   ---
#        reuse B5's current date-only format for the new B6 cell     ---
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("C5:G5").Copy()
$ws.Range("C6:G6").PasteSpecial(-4122)

# --- 4. Re-format B5 so it matches B2:B4 (datetime format) ---
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# --- 5. Fill in the new row's values ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 43987
$ws.Range("C6").Value = 110026
$ws.Range("D6").Value = 166049
$ws.Range("E6").Value = 48822
$ws.Range("F6").Value = 13170
$ws.Range("G6").Value = 34.03
